$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed (old rows 26 and 27);
# this shifts nothing else since they were the last two rows, and
# the data that was in row 27 is superseded by new literal values written to row 25 below.
$ws.Rows("26:27").Delete()

# Apply the individual cell value updates (rows 2-25)
$ws.Range("I2").Value = 2.2
$ws.Range("K2").Value = 2.2
$ws.Range("Q2").Value = 2.03
$ws.Range("R2").Value = 1.87
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 12
$ws.Range("AF2").Value = 41
$ws.Range("AL2").Value = 17
$ws.Range("AM2").Value = 26
$ws.Range("AU2").Value = 7.5
$ws.Range("J3").Value = 2.88
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3.2
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.91
$ws.Range("AC3").Value = 8
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 9.5
$ws.Range("AP3").Value = 23
$ws.Range("AT3").Value = 2.63
$ws.Range("AV3").Value = 51
$ws.Range("BD3").Value = 151
$ws.Range("G4").Value = 2.45
$ws.Range("I4").Value = 2.88
$ws.Range("J4").Value = 3.1
$ws.Range("W4").Value = 9
$ws.Range("AH4").Value = 9.5
$ws.Range("AL4").Value = 21
$ws.Range("G5").Value = 2.75
$ws.Range("I5").Value = 2.42
$ws.Range("J5").Value = 3.3
$ws.Range("L5").Value = 2.95
$ws.Range("P5").Value = 3.3
$ws.Range("S5").Value = 1.38
$ws.Range("T5").Value = 2.82
$ws.Range("W5").Value = 9.25
$ws.Range("Y5").Value = 10
$ws.Range("AB5").Value = 29
$ws.Range("AH5").Value = 8.75
$ws.Range("AL5").Value = 19
$ws.Range("AN5").Value = 4.75
$ws.Range("AT5").Value = 2.82
$ws.Range("AX5").Value = 12.5
$ws.Range("AY5").Value = 19.5
$ws.Range("AZ5").Value = 50
$ws.Range("BA5").Value = 80
$ws.Range("BD9").Value = 151
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 3.25
$ws.Range("J12").Value = 2.88
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.4
$ws.Range("X12").Value = 10
$ws.Range("AI12").Value = 17
$ws.Range("AL12").Value = 29
$ws.Range("AO12").Value = 12
$ws.Range("AX12").Value = 19
$ws.Range("AY12").Value = 29
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 9
$ws.Range("G18").Value = 2.6
$ws.Range("I18").Value = 2.63
$ws.Range("J18").Value = 3.2
$ws.Range("L18").Value = 3.25
$ws.Range("Q18").Value = 1.98
$ws.Range("R18").Value = 1.88
$ws.Range("W18").Value = 9
$ws.Range("AH18").Value = 9.5
$ws.Range("AQ18").Value = 41
$ws.Range("AZ18").Value = 51
$ws.Range("M20").Value = 1.08
$ws.Range("N20").Value = 8
$ws.Range("O20").Value = 1.4
$ws.Range("P20").Value = 2.75
$ws.Range("Q20").Value = 2.25
$ws.Range("R20").Value = 1.62
$ws.Range("G21").Value = 4.1
$ws.Range("J21").Value = 4.33
$ws.Range("L21").Value = 2.63
$ws.Range("N21").Value = 9.5
$ws.Range("W21").Value = 12
$ws.Range("X21").Value = 21
$ws.Range("AA21").Value = 34
$ws.Range("AK21").Value = 17
$ws.Range("G22").Value = 2.88
$ws.Range("I22").Value = 2.75
$ws.Range("J22").Value = 3.5
$ws.Range("M22").Value = 1.1
$ws.Range("N22").Value = 7
$ws.Range("U22").Value = 1.95
$ws.Range("V22").Value = 1.8
$ws.Range("AA22").Value = 26
$ws.Range("AB22").Value = 41
$ws.Range("AC22").Value = 7
$ws.Range("AG22").Value = 351
$ws.Range("AH22").Value = 7.5
$ws.Range("AI22").Value = 12
$ws.Range("AK22").Value = 26
$ws.Range("AN22").Value = 4.75
$ws.Range("AO22").Value = 17
$ws.Range("AS22").Value = 251
$ws.Range("BB22").Value = 251
$ws.Range("G24").Value = 2.38
$ws.Range("H24").Value = 3.2
$ws.Range("I24").Value = 3.1
$ws.Range("J24").Value = 3
$ws.Range("AH24").Value = 10
$ws.Range("AN24").Value = 4.33
$ws.Range("AY24").Value = 26
$ws.Range("A25").Value = 'bX7Spsu0'
$ws.Range("C25").Value = '16:30'
$ws.Range("D25").Value = 'SWITZERLAND - SUPER LEAGUE'
$ws.Range("E25").Value = 'Yverdon'
$ws.Range("F25").Value = 'Lausanne'
$ws.Range("G25").Value = 3.9
$ws.Range("H25").Value = 3.6
$ws.Range("I25").Value = 1.9
$ws.Range("J25").Value = 4
$ws.Range("K25").Value = 2.3
$ws.Range("L25").Value = 2.5
$ws.Range("M25").Value = 1.04
$ws.Range("N25").Value = 13
$ws.Range("O25").Value = 1.2
$ws.Range("P25").Value = 4.33
$ws.Range("Q25").Value = 1.7
$ws.Range("R25").Value = 2.1
$ws.Range("S25").Value = 1.33
$ws.Range("T25").Value = 3.25
$ws.Range("U25").Value = 1.58
$ws.Range("V25").Value = 2.2
$ws.Range("W25").Value = 13
$ws.Range("X25").Value = 21
$ws.Range("Y25").Value = 13
$ws.Range("Z25").Value = 41
$ws.Range("AA25").Value = 29
$ws.Range("AB25").Value = 34
$ws.Range("AC25").Value = 13
$ws.Range("AD25").Value = 7
$ws.Range("AE25").Value = 13
$ws.Range("AF25").Value = 41
$ws.Range("AG25").Value = 151
$ws.Range("AH25").Value = 9
$ws.Range("AI25").Value = 10
$ws.Range("AJ25").Value = 8.5
$ws.Range("AK25").Value = 17
$ws.Range("AL25").Value = 15
$ws.Range("AM25").Value = 21
$ws.Range("AN25").Value = 6
$ws.Range("AO25").Value = 19
$ws.Range("AP25").Value = 23
$ws.Range("AQ25").Value = 67
$ws.Range("AR25").Value = 67
$ws.Range("AS25").Value = 151
$ws.Range("AT25").Value = 3.25
$ws.Range("AU25").Value = 7.5
$ws.Range("AV25").Value = 41
$ws.Range("AW25").Value = 4
$ws.Range("AX25").Value = 10
$ws.Range("AY25").Value = 19
$ws.Range("AZ25").Value = 29
$ws.Range("BA25").Value = 41
$ws.Range("BB25").Value = 101
$ws.Range("BC25").Value = 81
$ws.Range("BD25").Value = 81
